# "Cleaned up and styled streamlit app"
#
# - Log: note that the AWS upload is confirmed working, with the external IP.
# - To Do: drop the two finished items ("Remove duplicates for koaatuu
#   locations" and "How to make the streamlit app always running" -- the
#   latter's slot is reused for a new item, "Read about TMUX"), widen the
#   notes column, and drop the now-obsolete hyperlink cell/formatting.

$wb = $excel.ActiveWorkbook

$wsLog = $wb.Worksheets.Item("Log")
$wsToDo = $wb.Worksheets.Item("To Do")

# --- Log sheet -------------------------------------------------------
$wsLog.Range("B7").Value = "Uploaded the app to the AWS EC2 and checked that it is working. The external IP address is http://18.117.238.29:8501/"

# --- To Do sheet -------------------------------------------------------
# "Remove duplicates for koaatuu locations" is done -- remove its row,
# shifting the remaining items up.
$wsToDo.Rows.Item(2).Delete()

# The old top item ("How to make the streamlit app always running") is
# done too; replace it with a freshly-added item.
$wsToDo.Range("A1").Value = "Read about TMUX"

# Drop the (now finished/irrelevant) hyperlink and its cell/formatting.
foreach ($link in $wsToDo.Hyperlinks) {
    $link.Delete()
}
$wsToDo.Range("B1").Clear()
$wsToDo.Rows.Item(1).AutoFit()

# Style: widen the notes column a bit.
$wsToDo.Columns.Item(2).ColumnWidth = 106.9

# --- restore cursor positions (To Do stays the active tab) -------------
$wsLog.Range("B16").Select() | Out-Null
$wsToDo.Range("A6").Select() | Out-Null
